# Generate Report for Handback
#
# A new handback-report generation run refreshed the timestamps for the
# "03bcaafe-f860-45ae-8450-eb80b0940753" file (the row that was just handed
# back / reprocessed), while the "17b0988e-966e-4b87-b1bf-6b03ee74d6cc" file
# row keeps its previous values.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-27 22:48:27"

# --- zh-cn sheet: Correspond Handoff / Handback DateTime (columns H & K) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-27 22:48:22"
$zhcn.Range("K2").Value = "2016-08-27 22:48:39"

# --- de-de sheet: Correspond Handoff / Handback DateTime (columns H & K) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-27 22:48:27"
$dede.Range("K2").Value = "2016-08-27 22:48:45"
